# Apply odds-data updates from the FlashScore 2025-05-20 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 2.75
$ws.Range("P2").Value = 1.36
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 1.67
$ws.Range("S2").Value = 2.1
$ws.Range("T2").Value = 10
$ws.Range("U2").Value = 15
$ws.Range("X2").Value = 21
$ws.Range("Y2").Value = 29
$ws.Range("Z2").Value = 11
$ws.Range("AD2").Value = 201
$ws.Range("AE2").Value = 11
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 29
$ws.Range("AI2").Value = 23

# Row 3
$ws.Range("I3").Value = 6
$ws.Range("V3").Value = 9.5
$ws.Range("W3").Value = 12
$ws.Range("AD3").Value = 151

# Row 4
$ws.Range("G4").Value = 2.9
$ws.Range("I4").Value = 2.35
$ws.Range("K4").Value = 10
$ws.Range("L4").Value = 1.3
$ws.Range("M4").Value = 3.4
$ws.Range("N4").Value = 2.03
$ws.Range("O4").Value = 1.78
$ws.Range("Z4").Value = 10
$ws.Range("AI4").Value = 19

# Row 8
$ws.Range("G8").Value = 6.5
$ws.Range("I8").Value = 1.5
$ws.Range("J8").Value = 1.07
$ws.Range("K8").Value = 9
$ws.Range("U8").Value = 34
$ws.Range("AG8").Value = 9.5

# Row 9
$ws.Range("G9").Value = 4.5
$ws.Range("H9").Value = 3.2
$ws.Range("I9").Value = 1.9
$ws.Range("N9").Value = 1.98
$ws.Range("O9").Value = 1.83
$ws.Range("P9").Value = 1.4
$ws.Range("Q9").Value = 2.75
$ws.Range("X9").Value = 34
$ws.Range("AA9").Value = 6
$ws.Range("AB9").Value = 13
$ws.Range("AF9").Value = 9
$ws.Range("AG9").Value = 9

# Row 16
$ws.Range("G16").Value = 1.87
$ws.Range("I16").Value = 3.95
$ws.Range("N16").Value = 2.05
$ws.Range("O16").Value = 1.6
$ws.Range("P16").Value = 1.39
$ws.Range("Q16").Value = 2.42
$ws.Range("T16").Value = 5.2
$ws.Range("U16").Value = 6.9
$ws.Range("V16").Value = 7.2
$ws.Range("W16").Value = 12.5
$ws.Range("X16").Value = 13.5
$ws.Range("Y16").Value = 25
$ws.Range("Z16").Value = 7.8
$ws.Range("AE16").Value = 8.5
$ws.Range("AG16").Value = 11
$ws.Range("AH16").Value = 45
$ws.Range("AI16").Value = 30

# Row 17
$ws.Range("G17").Value = 1.62
$ws.Range("H17").Value = 3.35
$ws.Range("I17").Value = 5.2
$ws.Range("N17").Value = 1.88
$ws.Range("T17").Value = 5.3
$ws.Range("U17").Value = 6.2
$ws.Range("W17").Value = 10
$ws.Range("Y17").Value = 21
$ws.Range("Z17").Value = 9
$ws.Range("AA17").Value = 5.8
$ws.Range("AC17").Value = 55
$ws.Range("AD17").Value = 350
$ws.Range("AE17").Value = 11.5
$ws.Range("AF17").Value = 26
$ws.Range("AG17").Value = 13.5
$ws.Range("AH17").Value = 75

# Row 24
$ws.Range("G24").Value = 1.2
$ws.Range("J24").Value = 21
$ws.Range("K24").Value = 1.03
$ws.Range("P24").Value = 1.2
$ws.Range("Q24").Value = 4.33
